# testDuplicate.xlsx fix-up: correct the sample data that shipped with a
# copy/paste bug (stale "Pablo Pinto" row baked in twice) and rename the
# email column header to match the rest of the loader's fixtures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -----------------------------------------------------
$ws.Range("B1").Value = "Email"

# --- Row 2 ------------------------------------------------------------
$ws.Range("A2").Value = "Pedro"
$ws.Range("B2").Value = "pedro@gmail.com"
$ws.Range("C2").Value = "56894512M"
$ws.Range("D2").Value = 12.568
$ws.Range("E2").Value = 15.26

# --- Row 3 (duplicate of row 2) ---------------------------------------
$ws.Range("A3").Value = "Pedro"
$ws.Range("B3").Value = "pedro@gmail.com"
$ws.Range("C3").Value = "56894512M"
$ws.Range("D3").Value = 12.568
$ws.Range("E3").Value = 15.26

# --- Selection: leave the cursor on A3, like the saved file -----------
$ws.Range("A3").Select()
